# Splits literal "<exp>...</exp>" markup inside run text into separate
# runs so the <exp> / </exp> tag delimiters get their own (Courier New,
# gray, small) character formatting, matching the "expan" tag styling
# already used elsewhere in this transcription for other inline tags
# (<ab>, <m>, <ms>, <page>, ...).

$d = $word.ActiveDocument

# Word's Font.Color (wdColor) is packed 0x00BBGGRR (BGR), not RGB, so
# convert the RGB hex values used elsewhere in the doc's palette.
function RgbToWdColor($hex) {
  $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
  $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
  $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
  return ($b * 65536) + ($g * 256) + $r
}

$expColor = RgbToWdColor "a9a9a9"
$delColor = RgbToWdColor "a91111"

# Finds literal "<exp>" / "</exp>" substrings inside $text (which starts
# at document offset $baseStart) and restyles just those tag delimiters,
# leaving the surrounding word-text runs with their original formatting.
function Format-ExpTags($doc, $baseStart, $text) {
  $openTag = "<exp>"
  $closeTag = "</exp>"
  $searchFrom = 0
  while ($true) {
    $openIdx = $text.IndexOf($openTag, $searchFrom)
    if ($openIdx -lt 0) { break }
    $closeIdx = $text.IndexOf($closeTag, $openIdx)
    if ($closeIdx -lt 0) { break }

    $openStart = $baseStart + $openIdx
    $openEnd = $openStart + $openTag.Length
    $openRange = $doc.Range($openStart, $openEnd)
    $openRange.Font.Name = "Courier New"
    $openRange.Font.Size = 7
    $openRange.Font.Color = $expColor

    $closeStart = $baseStart + $closeIdx
    $closeEnd = $closeStart + $closeTag.Length
    $closeRange = $doc.Range($closeStart, $closeEnd)
    $closeRange.Font.Name = "Courier New"
    $closeRange.Font.Size = 7
    $closeRange.Font.Color = $expColor

    $searchFrom = $closeIdx + $closeTag.Length
  }
}

# Locates $needle (a unique literal substring in the document) and runs
# Format-ExpTags across the whole matched range.
function Expand-ExpIn($doc, $needle) {
  $r = $doc.Content
  $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
  Format-ExpTags $doc $r.Start $r.Text
}

Expand-ExpIn $d "P<exp>rens</exp> <m>arsenic"
Expand-ExpIn $d " mediocrem<exp>ent</exp> subtille, 2 "
Expand-ExpIn $d "P<exp>rens</exp> de ces materiaulx, des meilleurs que tu pourras trouver. Poises"
Expand-ExpIn $d "les co<exp>mm</exp>e cy dessus est dict et les piles separement, vous tenant"
Expand-ExpIn $d ", pour fortifier v<exp>ost</exp>re "
Expand-ExpIn $d "En g<exp>e</exp>n<exp>er</exp>al,"
Expand-ExpIn $d "Le grain est co<exp>mm</exp>e"

# The lone "4" inside "<del>4</del>" gains the same small Courier New
# styling used for the other <del> markup in that run (color a91111),
# but at half size (sz/szCs 12 rather than 18).
$r = $d.Content
$r.Find.Execute("<del>4</del>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fourStart = $r.Start + 5
$fourEnd = $fourStart + 1
$fourRange = $d.Range($fourStart, $fourEnd)
$fourRange.Font.Name = "Courier New"
$fourRange.Font.Size = 6
$fourRange.Font.Color = $delColor
